$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6998.5
$ws.Range("I76").Value = 6998.5
$ws.Range("K76").Value = 6998.5
$ws.Range("M76").Value = -6683.5
$ws.Range("H79").Value = 6998.5
$ws.Range("I79").Value = 6998.5
$ws.Range("K79").Value = 6998.5
$ws.Range("M79").Value = -5906.5
$ws.Range("H113").Value = 4981.1113
$ws.Range("I113").Value = 2668.3333
$ws.Range("J113").Value = 6137.5
$ws.Range("K113").Value = 2668.3333
$ws.Range("L113").Value = 6137.5
$ws.Range("M113").Value = 585.6667000000002
$ws.Range("N113").Value = -12645.5
$ws.Range("H137").Value = 3523.7551
$ws.Range("J137").Value = 4956.607
$ws.Range("L137").Value = 14869.821
$ws.Range("N137").Value = -19969.821
$ws.Range("H138").Value = 4175.65
$ws.Range("J138").Value = 4609.737
$ws.Range("L138").Value = 13829.211
$ws.Range("N138").Value = -24109.211

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 45533628
$ws.Range("I45").Value = 107863.25
$ws.Range("J45").Value = 166669000
$ws.Range("K45").Value = 107863.25
$ws.Range("L45").Value = 166669000
$ws.Range("M45").Value = -107486.25
$ws.Range("N45").Value = -166669754
$ws.Range("H63").Value = 2270.5715
$ws.Range("I63").Value = 2270.5715
$ws.Range("K63").Value = 2270.5715
$ws.Range("M63").Value = -1584.5715
$ws.Range("H66").Value = 2270.5715
$ws.Range("I66").Value = 2270.5715
$ws.Range("K66").Value = 11352.8575
$ws.Range("M66").Value = -7920.8575
$ws.Range("H102").Value = 2029
$ws.Range("I102").Value = 1820.4242
$ws.Range("K102").Value = 1820.4242
$ws.Range("M102").Value = -198.4241999999999
$ws.Range("H110").Value = 2080.9524
$ws.Range("I110").Value = 2042.1052
$ws.Range("K110").Value = 2042.1052
$ws.Range("M110").Value = 2.894800000000032

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 21933210
$ws.Range("I20").Value = 26885576
$ws.Range("K20").Value = 26885576
$ws.Range("M20").Value = -26885329

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25001068
$ws.Range("I16").Value = 31251066
$ws.Range("K16").Value = 31251066
$ws.Range("M16").Value = -31250779
$ws.Range("H31").Value = 2624.4707
$ws.Range("I31").Value = 2887.25
$ws.Range("J31").Value = 2543.6155
$ws.Range("K31").Value = 2887.25
$ws.Range("L31").Value = 2543.6155
$ws.Range("M31").Value = -2592.25
$ws.Range("N31").Value = -3133.6155
$ws.Range("H34").Value = 2624.4707
$ws.Range("I34").Value = 2887.25
$ws.Range("J34").Value = 2543.6155
$ws.Range("K34").Value = 2887.25
$ws.Range("L34").Value = 2543.6155
$ws.Range("M34").Value = -2685.25
$ws.Range("N34").Value = -2947.6155
$ws.Range("H58").Value = 2505.4211
$ws.Range("I58").Value = 2100
$ws.Range("J58").Value = 2692.5386
$ws.Range("K58").Value = 2100
$ws.Range("L58").Value = 2692.5386
$ws.Range("M58").Value = -1897
$ws.Range("N58").Value = -3098.5386
$ws.Range("H99").Value = 83336920
$ws.Range("J99").Value = 5999.1665
$ws.Range("L99").Value = 5999.1665
$ws.Range("N99").Value = -8995.166499999999
$ws.Range("H113").Value = 25001068
$ws.Range("I113").Value = 31251066
$ws.Range("K113").Value = 31251066
$ws.Range("M113").Value = -31248896
$ws.Range("H126").Value = 83336920
$ws.Range("J126").Value = 5999.1665
$ws.Range("L126").Value = 17997.4995
$ws.Range("N126").Value = -22937.4995
$ws.Range("H132").Value = 4219
$ws.Range("I132").Value = 4398.75
$ws.Range("K132").Value = 13196.25
$ws.Range("M132").Value = -10666.25
$ws.Range("H134").Value = 3227.372
$ws.Range("I134").Value = 2376.3333
$ws.Range("K134").Value = 7128.999899999999
$ws.Range("M134").Value = -4593.999899999999
$ws.Range("H136").Value = 2505.4211
$ws.Range("I136").Value = 2100
$ws.Range("J136").Value = 2692.5386
$ws.Range("K136").Value = 6300
$ws.Range("L136").Value = 8077.6158
$ws.Range("M136").Value = -3750
$ws.Range("N136").Value = -13177.6158

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1473860.9
$ws.Range("J68").Value = 1698331.8
$ws.Range("L68").Value = 5094995.4
$ws.Range("N68").Value = -5096617.4
$ws.Range("H71").Value = 1473860.9
$ws.Range("J71").Value = 1698331.8
$ws.Range("L71").Value = 15284986.2
$ws.Range("N71").Value = -15293098.2
$ws.Range("H92").Value = 1959.1428
$ws.Range("J92").Value = 1959.1428
$ws.Range("L92").Value = 5877.428400000001
$ws.Range("N92").Value = -8373.428400000001
$ws.Range("H107").Value = 4665.926
$ws.Range("J107").Value = 8379.143
$ws.Range("L107").Value = 25137.429
$ws.Range("N107").Value = -28977.429
$ws.Range("H114").Value = 3111.8
$ws.Range("I114").Value = 2387.6
$ws.Range("J114").Value = 3836
$ws.Range("K114").Value = 7162.799999999999
$ws.Range("L114").Value = 11508
$ws.Range("M114").Value = -3908.799999999999
$ws.Range("N114").Value = -18016
$ws.Range("H130").Value = 8468
$ws.Range("I130").Value = 8468
$ws.Range("K130").Value = 25404
$ws.Range("M130").Value = -20384
$ws.Range("H132").Value = 6284.206
$ws.Range("I132").Value = 5116.3335
$ws.Range("K132").Value = 46047.0015
$ws.Range("M132").Value = -43517.0015

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 28578874
$ws.Range("I102").Value = 33337778
$ws.Range("K102").Value = 33337778
$ws.Range("M102").Value = -33336156
$ws.Range("H132").Value = 4500
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -13060.0001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4021.7144
$ws.Range("I7").Value = 3829.4
$ws.Range("K7").Value = 3829.4
$ws.Range("M7").Value = -3717.4
$ws.Range("H40").Value = 41535.156
$ws.Range("I40").Value = 59860.81
$ws.Range("K40").Value = 59860.81
$ws.Range("M40").Value = -59724.81
$ws.Range("H55").Value = 409.27274
$ws.Range("I55").Value = 337.2143
$ws.Range("K55").Value = 337.2143
$ws.Range("M55").Value = -164.2143
$ws.Range("H100").Value = 2087.8333
$ws.Range("I100").Value = 2615
$ws.Range("K100").Value = 2615
$ws.Range("M100").Value = -2074
$ws.Range("H122").Value = 5529.8
$ws.Range("I122").Value = 5162.25
$ws.Range("K122").Value = 15486.75
$ws.Range("M122").Value = -13036.75
$ws.Range("H126").Value = 4021.7144
$ws.Range("I126").Value = 3829.4
$ws.Range("K126").Value = 11488.2
$ws.Range("M126").Value = -9018.200000000001
$ws.Range("H132").Value = 10562.963
$ws.Range("I132").Value = 12784.5
$ws.Range("J132").Value = 8170.5386
$ws.Range("K132").Value = 38353.5
$ws.Range("L132").Value = 24511.6158
$ws.Range("M132").Value = -35823.5
$ws.Range("N132").Value = -29571.6158
$ws.Range("H136").Value = 6591.773
$ws.Range("I136").Value = 4639.3
$ws.Range("K136").Value = 13917.9
$ws.Range("M136").Value = -11367.9

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8666
$ws.Range("I62").Value = 10749
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 10749
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -10125
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 8666
$ws.Range("I65").Value = 10749
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 53745
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -50625
$ws.Range("N65").Value = -28740
$ws.Range("H81").Value = 5533.1665
$ws.Range("I81").Value = 5599.7144
$ws.Range("K81").Value = 11199.4288
$ws.Range("M81").Value = -10138.4288
$ws.Range("H84").Value = 5533.1665
$ws.Range("I84").Value = 5599.7144
$ws.Range("K84").Value = 55997.144
$ws.Range("M84").Value = -50693.144
$ws.Range("H107").Value = 684.2727
$ws.Range("I107").Value = 656.7778
$ws.Range("K107").Value = 1970.3334
$ws.Range("M107").Value = -50.33339999999998
$ws.Range("H122").Value = 10420160
$ws.Range("I122").Value = 3539.476
$ws.Range("K122").Value = 10618.428
$ws.Range("M122").Value = -8168.428
$ws.Range("H126").Value = 9989.429
$ws.Range("I126").Value = 11254.5
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 33763.5
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -31293.5
$ws.Range("N126").Value = -12137
$ws.Range("H132").Value = 6600
$ws.Range("I132").Value = 6723.077
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 20169.231
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -17639.231
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 20842264
$ws.Range("I136").Value = 23818062
$ws.Range("K136").Value = 71454186
$ws.Range("M136").Value = -71451636
